$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(13493,13493,13191,12237,11785,11785,10900,10617,10617,10522,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463,9463)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $values[$i]
}
